$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 10.972573550256008
$ws.Range("E5").Value = 1.6015398359937836
$ws.Range("H5").Value = 0.03469433409378291
$ws.Range("K5").Value = 7.258841774932164
$ws.Range("N5").Value = 1.378056117761644
$ws.Range("Q5").Value = 62.101350262295234
$ws.Range("T5").Value = 12.202172320009169
$ws.Range("W5").Value = 352.9755475969566
$ws.Range("Z5").Value = 144.67228942684875
$ws.Range("AC5").Value = 1.3206705324618147
$ws.Range("AF5").Value = 29.587931518346625
$ws.Range("AI5").Value = 8.319624656443828
$ws.Range("AL5").Value = 7.141975709761209
$ws.Range("AO5").Value = 1.3485377859735772
$ws.Range("AR5").Value = 36.21660983810309

$ws = $wb.Worksheets.Item(2)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 2.9948088737325573
$ws.Range("E5").Value = 1.7067612941289356
$ws.Range("H5").Value = 0.009320425660515096
$ws.Range("K5").Value = 4.957883013484282
$ws.Range("N5").Value = 0.8197687823681038
$ws.Range("Q5").Value = 22.96012398740883
$ws.Range("T5").Value = 3.115468191149978
$ws.Range("W5").Value = 71.99463558198715
$ws.Range("Z5").Value = 23.330549579910578
$ws.Range("AC5").Value = 1.2579262263692503
$ws.Range("AF5").Value = 30.72541414107303
$ws.Range("AI5").Value = 6.455700606913465
$ws.Range("AL5").Value = 5.002901968713873
$ws.Range("AO5").Value = 0.8298595040353839
$ws.Range("AR5").Value = 26.041973249452795

$ws = $wb.Worksheets.Item(3)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 14.65534651283338
$ws.Range("E5").Value = 1.8211441896428517
$ws.Range("H5").Value = 0.036225547688603595
$ws.Range("K5").Value = 6.484323118232638
$ws.Range("N5").Value = 1.5461927911090616
$ws.Range("Q5").Value = 102.69994848112457
$ws.Range("T5").Value = 16.509735078808315
$ws.Range("W5").Value = 384.21061677062517
$ws.Range("Z5").Value = 196.28743150359932
$ws.Range("AC5").Value = 1.6100749391445879
$ws.Range("AF5").Value = 26.92767964790576
$ws.Range("AI5").Value = 11.606858542424055
$ws.Range("AL5").Value = 6.3634336754175225
$ws.Range("AO5").Value = 1.5053963868571438
$ws.Range("AR5").Value = 36.73827256087299

$ws = $wb.Worksheets.Item(4)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 3.9409386573553182
$ws.Range("E5").Value = 1.2343736777814998
$ws.Range("H5").Value = 0.01856009593945782
$ws.Range("K5").Value = 4.038784947074429
$ws.Range("N5").Value = 0.8874576874207968
$ws.Range("Q5").Value = 26.47180005030692
$ws.Range("T5").Value = 4.074288751651805
$ws.Range("W5").Value = 74.30361297201556
$ws.Range("Z5").Value = 34.18530852284138
$ws.Range("AC5").Value = 1.4125548457344281
$ws.Range("AF5").Value = 23.619553042422794
$ws.Range("AI5").Value = 6.185863089543955
$ws.Range("AL5").Value = 3.992183724859212
$ws.Range("AO5").Value = 0.8834010277558461
$ws.Range("AR5").Value = 19.56328912886417

$ws = $wb.Worksheets.Item(5)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 1.5768725392625216
$ws.Range("E5").Value = 1.3757059512644896
$ws.Range("H5").Value = 0.008150080835435483
$ws.Range("K5").Value = 4.350165094306845
$ws.Range("N5").Value = 0.5877692031760554
$ws.Range("Q5").Value = 10.410163097824123
$ws.Range("T5").Value = 1.5865814777235199
$ws.Range("W5").Value = 61.61339453069405
$ws.Range("Z5").Value = 9.967361294816348
$ws.Range("AC5").Value = 1.4004212060550003
$ws.Range("AF5").Value = 24.305712439668664
$ws.Range("AI5").Value = 5.710155947066819
$ws.Range("AL5").Value = 4.374495327689856
$ws.Range("AO5").Value = 0.5945697092352988
$ws.Range("AR5").Value = 25.143732511211468

$ws = $wb.Worksheets.Item(6)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 7.197761418003572
$ws.Range("E5").Value = 1.488275012122904
$ws.Range("H5").Value = 0.020868473642696048
$ws.Range("K5").Value = 3.269768470638017
$ws.Range("N5").Value = 0.6325859244904279
$ws.Range("Q5").Value = 51.495632773319045
$ws.Range("T5").Value = 7.4749380743328855
$ws.Range("W5").Value = 171.8221244656365
$ws.Range("Z5").Value = 58.448212303081185
$ws.Range("AC5").Value = 1.1652797747008812
$ws.Range("AF5").Value = 15.637522798571107
$ws.Range("AI5").Value = 5.442384664833746
$ws.Range("AL5").Value = 3.1650240551534083
$ws.Range("AO5").Value = 0.6064677785738021
$ws.Range("AR5").Value = 14.764174896877414

$ws = $wb.Worksheets.Item(7)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.802036174556092
$ws.Range("E5").Value = 1.19029927249144
$ws.Range("H5").Value = 0.009885687882693991
$ws.Range("K5").Value = 4.260564380571885
$ws.Range("N5").Value = 0.44677350287988776
$ws.Range("Q5").Value = 4.817596680175335
$ws.Range("T5").Value = 0.8566340914159872
$ws.Range("W5").Value = 14.34123759109778
$ws.Range("Z5").Value = 2.910358564573255
$ws.Range("AC5").Value = 1.2402255360137697
$ws.Range("AF5").Value = 20.672034776696886
$ws.Range("AI5").Value = 3.3964491115390687
$ws.Range("AL5").Value = 4.285354366925018
$ws.Range("AO5").Value = 0.4504224601632172
$ws.Range("AR5").Value = 13.637772712034755

$ws = $wb.Worksheets.Item(8)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 2.475634856682594
$ws.Range("E5").Value = 1.4013956185309533
$ws.Range("H5").Value = 0.008123645060435907
$ws.Range("K5").Value = 3.1324707753173473
$ws.Range("N5").Value = 0.8560058481639121
$ws.Range("Q5").Value = 20.936078654869007
$ws.Range("T5").Value = 2.5252464310947533
$ws.Range("W5").Value = 40.61200004481879
$ws.Range("Z5").Value = 15.01890441229636
$ws.Range("AC5").Value = 1.3899107971593894
$ws.Range("AF5").Value = 20.502758292063685
$ws.Range("AI5").Value = 7.465574227800372
$ws.Range("AL5").Value = 3.1111119200109623
$ws.Range("AO5").Value = 0.8575323481663704
$ws.Range("AR5").Value = 16.546900290034447

$ws = $wb.Worksheets.Item(9)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.6557633019375756
$ws.Range("E5").Value = 1.4315099842358663
$ws.Range("H5").Value = 0.003896751571969302
$ws.Range("K5").Value = 2.477592907268719
$ws.Range("N5").Value = 0.4495149536973573
$ws.Range("Q5").Value = 5.647478585937864
$ws.Range("T5").Value = 0.6625964788711312
$ws.Range("W5").Value = 12.419575327620331
$ws.Range("Z5").Value = 2.2895402293496145
$ws.Range("AC5").Value = 1.5261002472489282
$ws.Range("AF5").Value = 18.623955500891977
$ws.Range("AI5").Value = 5.116075691367994
$ws.Range("AL5").Value = 2.484989880470473
$ws.Range("AO5").Value = 0.4525470484413791
$ws.Range("AR5").Value = 12.29093131401718

$ws = $wb.Worksheets.Item(10)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.36631300783081516
$ws.Range("E5").Value = 1.4704464007563731
$ws.Range("H5").Value = 0.0035065365860613503
$ws.Range("K5").Value = 2.0906664536597552
$ws.Range("N5").Value = 0.25681337985547115
$ws.Range("Q5").Value = 3.322152369737527
$ws.Range("T5").Value = 0.38697843609459087
$ws.Range("W5").Value = 5.54099586690226
$ws.Range("Z5").Value = 0.9898013960193915
$ws.Range("AC5").Value = 1.410768043635887
$ws.Range("AF5").Value = 13.506076590200893
$ws.Range("AI5").Value = 3.1678895187395493
$ws.Range("AL5").Value = 2.099091590521827
$ws.Range("AO5").Value = 0.2577570503869008
$ws.Range("AR5").Value = 6.742468235169849

$ws = $wb.Worksheets.Item(11)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.8450943500566424
$ws.Range("E5").Value = 1.463294940007196
$ws.Range("H5").Value = 0.004672330635982965
$ws.Range("K5").Value = 2.9246262246366075
$ws.Range("N5").Value = 0.5222914926647341
$ws.Range("Q5").Value = 6.159305968199315
$ws.Range("T5").Value = 0.8547834468972618
$ws.Range("W5").Value = 13.259418515287665
$ws.Range("Z5").Value = 3.8627184760759845
$ws.Range("AC5").Value = 1.6464267223195101
$ws.Range("AF5").Value = 21.76438568758553
$ws.Range("AI5").Value = 5.002880668940173
$ws.Range("AL5").Value = 2.9444451467495534
$ws.Range("AO5").Value = 0.5257854836584469
$ws.Range("AR5").Value = 12.008398004250951

$ws = $wb.Worksheets.Item(12)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 1.170787162589821
$ws.Range("E5").Value = 1.442069369626389
$ws.Range("H5").Value = 0.007530355716838267
$ws.Range("K5").Value = 2.8949869550278042
$ws.Range("N5").Value = 0.5868098108661427
$ws.Range("Q5").Value = 6.923002277631459
$ws.Range("T5").Value = 1.1612864879169902
$ws.Range("W5").Value = 12.344551092030244
$ws.Range("Z5").Value = 6.232245025592803
$ws.Range("AC5").Value = 1.0814392104167865
$ws.Range("AF5").Value = 18.216143528563816
$ws.Range("AI5").Value = 4.156229285055912
$ws.Range("AL5").Value = 2.9076317181960536
$ws.Range("AO5").Value = 0.5933370745056575
$ws.Range("AR5").Value = 7.984985507309852
